# ----------------------------------------------------------------------------
# Scheduled-runner refresh: re-pull current Universalis market prices for the
# Zalera-server Leve-profit workbook and recompute the dependent profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) on every sheet
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Values only - no formulas live in
# this workbook, so each affected cell is written directly.
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 286.72726
$ws.Cells.Item(8, 8).Value = 551.5
$ws.Cells.Item(8, 9).Value = 551.5
$ws.Cells.Item(8, 11).Value = 1654.5
$ws.Cells.Item(8, 13).Value = -1515.5
$ws.Cells.Item(11, 8).Value = 5354.778
$ws.Cells.Item(11, 9).Value = 5354.778
$ws.Cells.Item(11, 11).Value = 5354.778
$ws.Cells.Item(11, 13).Value = -5214.778
$ws.Cells.Item(19, 8).Value = 1202.6364
$ws.Cells.Item(19, 9).Value = 1272.7778
$ws.Cells.Item(19, 10).Value = 887
$ws.Cells.Item(19, 11).Value = 1272.7778
$ws.Cells.Item(19, 12).Value = 887
$ws.Cells.Item(19, 13).Value = -1097.7778
$ws.Cells.Item(19, 14).Value = -1237
$ws.Cells.Item(38, 8).Value = 562.5
$ws.Cells.Item(38, 10).Value = 4500
$ws.Cells.Item(38, 12).Value = 13500
$ws.Cells.Item(38, 14).Value = -14244
$ws.Cells.Item(43, 8).Value = 16856
$ws.Cells.Item(43, 10).Value = 9799.6
$ws.Cells.Item(43, 12).Value = 9799.6
$ws.Cells.Item(43, 14).Value = -9937.6
$ws.Cells.Item(86, 8).Value = 2249.7856
$ws.Cells.Item(86, 9).Value = 2212.375
$ws.Cells.Item(86, 10).Value = 2299.6667
$ws.Cells.Item(86, 11).Value = 2212.375
$ws.Cells.Item(86, 12).Value = 2299.6667
$ws.Cells.Item(86, 13).Value = -1089.375
$ws.Cells.Item(86, 14).Value = -4545.6667
$ws.Cells.Item(89, 8).Value = 2249.7856
$ws.Cells.Item(89, 9).Value = 2212.375
$ws.Cells.Item(89, 10).Value = 2299.6667
$ws.Cells.Item(89, 11).Value = 11061.875
$ws.Cells.Item(89, 12).Value = 11498.3335
$ws.Cells.Item(89, 13).Value = -5445.875
$ws.Cells.Item(89, 14).Value = -22730.3335
$ws.Cells.Item(99, 8).Value = 336
$ws.Cells.Item(99, 9).Value = 350
$ws.Cells.Item(99, 10).Value = 331.33334
$ws.Cells.Item(99, 11).Value = 1050
$ws.Cells.Item(99, 12).Value = 994.0000200000001
$ws.Cells.Item(99, 13).Value = 448
$ws.Cells.Item(99, 14).Value = -3990.00002
$ws.Cells.Item(132, 8).Value = 902.07275
$ws.Cells.Item(132, 9).Value = 763.4897999999999
$ws.Cells.Item(132, 11).Value = 2290.4694
$ws.Cells.Item(132, 13).Value = 239.5306
$ws.Cells.Item(137, 8).Value = 13165878
$ws.Cells.Item(137, 9).Value = 27778750
$ws.Cells.Item(137, 10).Value = 14292.6
$ws.Cells.Item(137, 11).Value = 83336250
$ws.Cells.Item(137, 12).Value = 42877.8
$ws.Cells.Item(137, 13).Value = -83333700
$ws.Cells.Item(137, 14).Value = -47977.8
$ws.Cells.Item(141, 8).Value = 2105.8333
$ws.Cells.Item(141, 9).Value = 1527
$ws.Cells.Item(141, 11).Value = 4581
$ws.Cells.Item(141, 13).Value = 599

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 13144774
$ws.Cells.Item(2, 9).Value = 15335070
$ws.Cells.Item(2, 10).Value = 3000
$ws.Cells.Item(2, 11).Value = 15335070
$ws.Cells.Item(2, 12).Value = 3000
$ws.Cells.Item(2, 13).Value = -15334957
$ws.Cells.Item(2, 14).Value = -3226
$ws.Cells.Item(45, 8).Value = 2609.7856
$ws.Cells.Item(45, 9).Value = 2436.4167
$ws.Cells.Item(45, 10).Value = 3650
$ws.Cells.Item(45, 11).Value = 2436.4167
$ws.Cells.Item(45, 12).Value = 3650
$ws.Cells.Item(45, 13).Value = -2059.4167
$ws.Cells.Item(45, 14).Value = -4404
$ws.Cells.Item(50, 8).Value = 290.33334
$ws.Cells.Item(50, 9).Value = 231.66667
$ws.Cells.Item(50, 10).Value = 349
$ws.Cells.Item(50, 11).Value = 231.66667
$ws.Cells.Item(50, 12).Value = 349
$ws.Cells.Item(50, 13).Value = 482.33333
$ws.Cells.Item(50, 14).Value = -1777
$ws.Cells.Item(74, 8).Value = 258808.4
$ws.Cells.Item(74, 10).Value = 3574.5417
$ws.Cells.Item(74, 12).Value = 3574.5417
$ws.Cells.Item(74, 14).Value = -5322.5417
$ws.Cells.Item(77, 8).Value = 258808.4
$ws.Cells.Item(77, 10).Value = 3574.5417
$ws.Cells.Item(77, 12).Value = 17872.7085
$ws.Cells.Item(77, 14).Value = -26608.7085
$ws.Cells.Item(116, 8).Value = 13144774
$ws.Cells.Item(116, 9).Value = 15335070
$ws.Cells.Item(116, 10).Value = 3000
$ws.Cells.Item(116, 11).Value = 15335070
$ws.Cells.Item(116, 12).Value = 3000
$ws.Cells.Item(116, 13).Value = -15332776
$ws.Cells.Item(116, 14).Value = -7588
$ws.Cells.Item(122, 8).Value = 1440.5416
$ws.Cells.Item(122, 9).Value = 1171.8182
$ws.Cells.Item(122, 11).Value = 3515.4546
$ws.Cells.Item(122, 13).Value = -1065.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 13144774
$ws.Cells.Item(3, 9).Value = 15335070
$ws.Cells.Item(3, 10).Value = 3000
$ws.Cells.Item(3, 11).Value = 15335070
$ws.Cells.Item(3, 12).Value = 3000
$ws.Cells.Item(3, 13).Value = -15334956
$ws.Cells.Item(3, 14).Value = -3228
$ws.Cells.Item(20, 8).Value = 2351.5908
$ws.Cells.Item(20, 9).Value = 2144
$ws.Cells.Item(20, 11).Value = 2144
$ws.Cells.Item(20, 13).Value = -1897
$ws.Cells.Item(99, 8).Value = 2200.2632
$ws.Cells.Item(99, 9).Value = 1733.75
$ws.Cells.Item(99, 11).Value = 1733.75
$ws.Cells.Item(99, 13).Value = -235.75
$ws.Cells.Item(107, 8).Value = 2018.1
$ws.Cells.Item(107, 9).Value = 1884.0834
$ws.Cells.Item(107, 11).Value = 1884.0834
$ws.Cells.Item(107, 13).Value = 35.91660000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 405.23077
$ws.Cells.Item(7, 9).Value = 340.66666
$ws.Cells.Item(7, 10).Value = 550.5
$ws.Cells.Item(7, 11).Value = 340.66666
$ws.Cells.Item(7, 12).Value = 550.5
$ws.Cells.Item(7, 13).Value = -227.66666
$ws.Cells.Item(7, 14).Value = -776.5
$ws.Cells.Item(31, 8).Value = 38465530
$ws.Cells.Item(31, 9).Value = 111113030
$ws.Cells.Item(31, 10).Value = 5084.8823
$ws.Cells.Item(31, 11).Value = 111113030
$ws.Cells.Item(31, 12).Value = 5084.8823
$ws.Cells.Item(31, 13).Value = -111112735
$ws.Cells.Item(31, 14).Value = -5674.8823
$ws.Cells.Item(34, 8).Value = 38465530
$ws.Cells.Item(34, 9).Value = 111113030
$ws.Cells.Item(34, 10).Value = 5084.8823
$ws.Cells.Item(34, 11).Value = 111113030
$ws.Cells.Item(34, 12).Value = 5084.8823
$ws.Cells.Item(34, 13).Value = -111112828
$ws.Cells.Item(34, 14).Value = -5488.8823
$ws.Cells.Item(62, 8).Value = 27997.6
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(65, 8).Value = 27997.6
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(122, 8).Value = 72417.71000000001
$ws.Cells.Item(122, 9).Value = 84320.75
$ws.Cells.Item(122, 11).Value = 252962.25
$ws.Cells.Item(122, 13).Value = -250512.25
$ws.Cells.Item(134, 8).Value = 7626.107
$ws.Cells.Item(134, 9).Value = 7811.1904
$ws.Cells.Item(134, 11).Value = 23433.5712
$ws.Cells.Item(134, 13).Value = -20898.5712
$ws.Cells.Item(140, 8).Value = 120000
$ws.Cells.Item(140, 10).Value = 120000
$ws.Cells.Item(140, 12).Value = 120000
$ws.Cells.Item(140, 14).Value = -130360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 2582.6428
$ws.Cells.Item(2, 9).Value = 102.22222
$ws.Cells.Item(2, 10).Value = 3757.5789
$ws.Cells.Item(2, 11).Value = 613.33332
$ws.Cells.Item(2, 12).Value = 22545.4734
$ws.Cells.Item(2, 13).Value = -500.33332
$ws.Cells.Item(2, 14).Value = -22771.4734
$ws.Cells.Item(12, 8).Value = 1100.16
$ws.Cells.Item(12, 10).Value = 984.4
$ws.Cells.Item(12, 12).Value = 2953.2
$ws.Cells.Item(12, 14).Value = -3299.2
$ws.Cells.Item(50, 8).Value = 1363.8
$ws.Cells.Item(50, 10).Value = 2999.25
$ws.Cells.Item(50, 12).Value = 8997.75
$ws.Cells.Item(50, 14).Value = -9959.75
$ws.Cells.Item(53, 8).Value = 1363.8
$ws.Cells.Item(53, 10).Value = 2999.25
$ws.Cells.Item(53, 12).Value = 8997.75
$ws.Cells.Item(53, 14).Value = -9959.75
$ws.Cells.Item(75, 8).Value = 1566.25
$ws.Cells.Item(75, 10).Value = 1802
$ws.Cells.Item(75, 12).Value = 5406
$ws.Cells.Item(75, 14).Value = -7402
$ws.Cells.Item(78, 8).Value = 1566.25
$ws.Cells.Item(78, 10).Value = 1802
$ws.Cells.Item(78, 12).Value = 16218
$ws.Cells.Item(78, 14).Value = -26202

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 59548.293
$ws.Cells.Item(2, 9).Value = 852.0769
$ws.Cells.Item(2, 11).Value = 852.0769
$ws.Cells.Item(2, 13).Value = -739.0769
$ws.Cells.Item(113, 8).Value = 23125
$ws.Cells.Item(113, 9).Value = 1589.625
$ws.Cells.Item(113, 10).Value = 44660.375
$ws.Cells.Item(113, 11).Value = 1589.625
$ws.Cells.Item(113, 12).Value = 44660.375
$ws.Cells.Item(113, 13).Value = 580.375
$ws.Cells.Item(113, 14).Value = -49000.375
$ws.Cells.Item(126, 8).Value = 2534.6667
$ws.Cells.Item(126, 9).Value = 2226.2856
$ws.Cells.Item(126, 11).Value = 6678.8568
$ws.Cells.Item(126, 13).Value = -4208.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 5811.7856
$ws.Cells.Item(93, 10).Value = 5269.5454
$ws.Cells.Item(93, 12).Value = 5269.5454
$ws.Cells.Item(93, 14).Value = -7765.5454
$ws.Cells.Item(122, 8).Value = 5150.154
$ws.Cells.Item(122, 9).Value = 4017.625
$ws.Cells.Item(122, 10).Value = 6962.2
$ws.Cells.Item(122, 11).Value = 12052.875
$ws.Cells.Item(122, 12).Value = 20886.6
$ws.Cells.Item(122, 13).Value = -9602.875
$ws.Cells.Item(122, 14).Value = -25786.6
$ws.Cells.Item(132, 8).Value = 5804.952
$ws.Cells.Item(132, 9).Value = 4802.6
$ws.Cells.Item(132, 10).Value = 8310.833000000001
$ws.Cells.Item(132, 11).Value = 14407.8
$ws.Cells.Item(132, 12).Value = 24932.499
$ws.Cells.Item(132, 13).Value = -11877.8
$ws.Cells.Item(132, 14).Value = -29992.499

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 7295.273
$ws.Cells.Item(81, 9).Value = 2928.4285
$ws.Cells.Item(81, 10).Value = 14937.25
$ws.Cells.Item(81, 11).Value = 5856.857
$ws.Cells.Item(81, 12).Value = 29874.5
$ws.Cells.Item(81, 13).Value = -4795.857
$ws.Cells.Item(81, 14).Value = -31996.5
$ws.Cells.Item(84, 8).Value = 7295.273
$ws.Cells.Item(84, 9).Value = 2928.4285
$ws.Cells.Item(84, 10).Value = 14937.25
$ws.Cells.Item(84, 11).Value = 29284.285
$ws.Cells.Item(84, 12).Value = 149372.5
$ws.Cells.Item(84, 13).Value = -23980.285
$ws.Cells.Item(84, 14).Value = -159980.5
$ws.Cells.Item(100, 8).Value = 850.04346
$ws.Cells.Item(100, 9).Value = 750.9474
$ws.Cells.Item(100, 11).Value = 1501.8948
$ws.Cells.Item(100, 13).Value = -960.8948
$ws.Cells.Item(122, 8).Value = 4893.909
$ws.Cells.Item(122, 9).Value = 4982.8
$ws.Cells.Item(122, 10).Value = 4005
$ws.Cells.Item(122, 11).Value = 14948.4
$ws.Cells.Item(122, 12).Value = 12015
$ws.Cells.Item(122, 13).Value = -12498.4
$ws.Cells.Item(122, 14).Value = -16915
$ws.Cells.Item(132, 8).Value = 6203.4443
$ws.Cells.Item(132, 9).Value = 4097.4287
$ws.Cells.Item(132, 10).Value = 6940.55
$ws.Cells.Item(132, 11).Value = 12292.2861
$ws.Cells.Item(132, 12).Value = 20821.65
$ws.Cells.Item(132, 13).Value = -9762.286100000001
$ws.Cells.Item(132, 14).Value = -25881.65
$ws.Cells.Item(133, 8).Value = 73266
$ws.Cells.Item(133, 10).Value = 73266
$ws.Cells.Item(133, 12).Value = 73266
$ws.Cells.Item(133, 14).Value = -83386
